$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 ("Custom domain name") is now complete: copy the formatting used by
# other completed rows (e.g. row 4) and fill in the resolved date.
$ws.Range("B4:D4").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)
$ws.Range("E4:F4").Copy()
$ws.Range("E17:F17").PasteSpecial(-4122)
$ws.Range("F17").Value = 44331

# Rename the "Error Warnings" task (row 22) to be more specific.
$ws.Range("B22").Value = "Error Warnings (Accessibility) "

# Update the active selection shown when the file was last saved.
[void]$ws.Range("H17").Select()
